$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 13
$ws.Range("D2").Value = 45096
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 23000
$ws.Range("P2").Value = 23000
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 1278
$ws.Range("T2").Value = 18

# Row 3 <- original row 17
$ws.Range("D3").Value = 45041
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("R3").Value = 'Región Metropolitana'
$ws.Range("S3").Value = 833
$ws.Range("T3").Value = 18

# Row 4 <- original row 18
$ws.Range("D4").Value = 45028
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18

# Row 5 <- original row 2
$ws.Range("D5").Value = 45037
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 889
$ws.Range("T5").Value = 18

# Row 6 <- original row 15
$ws.Range("D6").Value = 45033
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18

# Row 7 <- original row 8
$ws.Range("D7").Value = 44999
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 17500
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 972
$ws.Range("T7").Value = 18

# Row 8 <- original row 22
$ws.Range("D8").Value = 45043
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 833
$ws.Range("T8").Value = 18

# Row 9 <- original row 10
$ws.Range("D9").Value = 45021
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15500
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Provincia de Los Andes'
$ws.Range("S9").Value = 861
$ws.Range("T9").Value = 18

# Row 10 <- original row 19
$ws.Range("D10").Value = 45091
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 22000
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 1222
$ws.Range("T10").Value = 18

# Row 11 <- original row 4
$ws.Range("D11").Value = 45036
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 15500
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 861
$ws.Range("T11").Value = 18

# Row 12 <- original row 5
$ws.Range("D12").Value = 45020
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/caja 16 kilos'
$ws.Range("R12").Value = 'Provincia de Los Andes'
$ws.Range("S12").Value = 938
$ws.Range("T12").Value = 16

# Row 13 <- original row 6
$ws.Range("D13").Value = 45044
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 833
$ws.Range("T13").Value = 18

# Row 14 <- original row 3
$ws.Range("D14").Value = 45001
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17500
$ws.Range("Q14").Value = '$/caja 18 kilos'
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 972
$ws.Range("T14").Value = 18

# Row 15 <- original row 12
$ws.Range("D15").Value = 45050
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("Q15").Value = '$/caja 18 kilos'
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("S15").Value = 778
$ws.Range("T15").Value = 18

# Row 16 <- original row 14
$ws.Range("D16").Value = 45030
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 1000
$ws.Range("T16").Value = 18

# Row 17 <- original row 9
$ws.Range("D17").Value = 45049
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 833
$ws.Range("T17").Value = 18

# Row 18 <- original row 7
$ws.Range("D18").Value = 45002
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = 18000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 18000
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 1000
$ws.Range("T18").Value = 18

# Row 19 <- original row 20
$ws.Range("D19").Value = 45014
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = '$/caja 18 kilos'
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18

# Row 20 <- original row 11
$ws.Range("D20").Value = 45062
$ws.Range("M20").Value = 90
$ws.Range("N20").Value = 13000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 13444
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 747
$ws.Range("T20").Value = 18

# Row 22 <- original row 16
$ws.Range("D22").Value = 45099
$ws.Range("M22").Value = 40
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 22000
$ws.Range("Q22").Value = '$/caja 18 kilos'
$ws.Range("R22").Value = 'Región Metropolitana'
$ws.Range("S22").Value = 1222
$ws.Range("T22").Value = 18
